# Performance tab + AI
# Adds an "HFT_Profile" column (H) to the Sections sheet, tags every row
# with a profile (General / LowLatency / Network / Memory / Storage /
# Custom), fills in a few missing threshold values for the network rows,
# and turns on an AutoFilter on column H filtered down to "Storage".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New column header (H1) - style matches the other plain headers
#    (G1), i.e. same font/color as the rest of the non-bold header row.
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "HFT_Profile"
$ws.Range("H1").Font.Color = 0

# ---------------------------------------------------------------------
# 2) Per-row profile tag (column H) for every data row, rows 2-57.
# ---------------------------------------------------------------------
$profiles = @{
    2  = "General";    3  = "LowLatency"; 4  = "LowLatency"; 5  = "LowLatency";
    6  = "LowLatency"; 7  = "LowLatency"; 8  = "LowLatency"; 9  = "Network";
    10 = "Network";    11 = "Network";    12 = "General";    13 = "Memory";
    14 = "Memory";     15 = "Memory";     16 = "Memory";     17 = "General";
    18 = "LowLatency"; 19 = "Network";    20 = "LowLatency"; 21 = "General";
    22 = "General";    23 = "General";    24 = "General";    25 = "General";
    26 = "General";    27 = "General";    28 = "Network";    29 = "Network";
    30 = "Network";    31 = "Network";    32 = "Network";    33 = "Network";
    34 = "Network";    35 = "Network";    36 = "Network";    37 = "Network";
    38 = "Network";    39 = "Network";    40 = "Network";    41 = "Network";
    42 = "Network";    43 = "Network";    44 = "Network";    45 = "Network";
    46 = "Network";    47 = "Network";    48 = "Network";    49 = "Network";
    50 = "Network";    51 = "Network";    52 = "Network";    53 = "Network";
    54 = "Storage";    55 = "Storage";    56 = "Storage";    57 = "Custom"
}

foreach ($r in 2..57) {
    $ws.Cells.Item($r, 8).Value = $profiles[$r]
}

# ---------------------------------------------------------------------
# 3) Fill in a handful of newly-added threshold values on the Network
#    rows (Threshold_Min / Threshold_Max, columns E/F).
# ---------------------------------------------------------------------
$ws.Range("E35").Value = 100
$ws.Range("F35").Value = 120

$ws.Range("E36").Value = 4000000
$ws.Range("F36").Value = 6000000

$ws.Range("E37").Value = 4000000
$ws.Range("F37").Value = 6000000

$ws.Range("F38").Value = 5

# ---------------------------------------------------------------------
# 4) Turn on the AutoFilter across the (now) A1:H57 range and filter
#    column H (the 8th column) down to just "Storage".
# ---------------------------------------------------------------------
[void]$ws.Range("A1:H57").AutoFilter(8, @("Storage"), 7)

# ---------------------------------------------------------------------
# 5) Hidden worksheet-scoped defined name Excel writes for the active
#    filter range.
# ---------------------------------------------------------------------
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sections!`$A`$1:`$H`$57")
$fdb.Visible = $false

# ---------------------------------------------------------------------
# 6) Restore the view: select H1 (matches the selection left behind by
#    the author after adding/filtering the column).
# ---------------------------------------------------------------------
[void]$ws.Range("H1").Select()
